$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data-driven cell updates reproducing the upstream "cryptos list" refresh.
# Each entry is a cell reference + its new text value. We force the cell to
# Text format before writing so numeric-looking strings (e.g. "239.06") are
# NOT auto-converted to numbers by Excel (matching the source file, where
# every data cell is stored as text/inlineStr). We reset the style back to
# "Normal" afterwards so no stray number-format styling is left on the cell.
$changes = @(
    @{ Cell = "D2"; Value = "29.326.27" }
    @{ Cell = "E2"; Value = "  -0.19%  " }
    @{ Cell = "D3"; Value = "1.841.11" }
    @{ Cell = "E3"; Value = "  -0.36%  " }
    @{ Cell = "D4"; Value = "0.9988" }
    @{ Cell = "D5"; Value = "239.06" }
    @{ Cell = "E5"; Value = "  -0.64%  " }
    @{ Cell = "D6"; Value = "0.6278" }
    @{ Cell = "E6"; Value = "  -0.59%  " }
    @{ Cell = "D8"; Value = "0.07517" }
    @{ Cell = "E8"; Value = "  -0.75%  " }
    @{ Cell = "D9"; Value = "0.2931" }
    @{ Cell = "D10"; Value = "24.34" }
    @{ Cell = "E10"; Value = "  -0.60%  " }
    @{ Cell = "D11"; Value = "0.07688" }
    @{ Cell = "E11"; Value = "  -0.49%  " }
    @{ Cell = "D12"; Value = "1.828.41" }
    @{ Cell = "E12"; Value = "  -7.89%  " }
    @{ Cell = "D13"; Value = "4.979" }
    @{ Cell = "E13"; Value = "  -0.12%  " }
    @{ Cell = "D14"; Value = "0.6771" }
    @{ Cell = "E14"; Value = "  -1.12%  " }
    @{ Cell = "E15"; Value = "  +4.82%  " }
    @{ Cell = "D16"; Value = "82.78" }
    @{ Cell = "E16"; Value = "  -0.05%  " }
    @{ Cell = "D17"; Value = "2.079.62" }
    @{ Cell = "E17"; Value = "  -8.16%  " }
    @{ Cell = "D18"; Value = "6.096" }
    @{ Cell = "E18"; Value = "  -1.49%  " }
    @{ Cell = "D19"; Value = "29.356.71" }
    @{ Cell = "E19"; Value = "  -0.31%  " }
    @{ Cell = "D20"; Value = "227.26" }
    @{ Cell = "E20"; Value = "  -1.67%  " }
    @{ Cell = "E21"; Value = "  -0.91%  " }
    @{ Cell = "D22"; Value = "1.000" }
    @{ Cell = "E22"; Value = "  +0.01%  " }
    @{ Cell = "D23"; Value = "7.403" }
    @{ Cell = "E23"; Value = "  -2.51%  " }
    @{ Cell = "E24"; Value = "  +0.03%  " }
    @{ Cell = "D25"; Value = "156.59" }
    @{ Cell = "E25"; Value = "  +1.18%  " }
    @{ Cell = "D26"; Value = "0.1383" }
    @{ Cell = "E26"; Value = "  -0.86%  " }
    @{ Cell = "E27"; Value = "  -1.12%  " }
    @{ Cell = "D28"; Value = "17.57" }
    @{ Cell = "E28"; Value = "  -0.48%  " }
    @{ Cell = "D29"; Value = "1.459" }
    @{ Cell = "E29"; Value = "  -0.82%  " }
    @{ Cell = "D30"; Value = "1.276" }
    @{ Cell = "E30"; Value = "  +0.72%  " }
    @{ Cell = "D31"; Value = "0.05614" }
    @{ Cell = "E31"; Value = "  -3.36%  " }
    @{ Cell = "E32"; Value = "  -0.66%  " }
    @{ Cell = "D33"; Value = "4.014" }
    @{ Cell = "E33"; Value = "  -0.11%  " }
    @{ Cell = "E34"; Value = "  -2.12%  " }
    @{ Cell = "E35"; Value = "  -0.47%  " }
    @{ Cell = "D36"; Value = "0.7069" }
    @{ Cell = "E37"; Value = "  -0.26%  " }
    @{ Cell = "D38"; Value = "1.238.84" }
    @{ Cell = "E38"; Value = "  -0.85%  " }
    @{ Cell = "D39"; Value = "0.01802" }
    @{ Cell = "E39"; Value = "  -0.24%  " }
    @{ Cell = "D40"; Value = "2.757" }
    @{ Cell = "E40"; Value = "  -1.31%  " }
    @{ Cell = "D41"; Value = "6.246" }
    @{ Cell = "E41"; Value = "  +2.51%  " }
    @{ Cell = "D42"; Value = "0.9010" }
    @{ Cell = "E42"; Value = "  -0.47%  " }
    @{ Cell = "D43"; Value = "0.9991" }
    @{ Cell = "E43"; Value = "  -0.03%  " }
    @{ Cell = "D44"; Value = "101.77" }
    @{ Cell = "E44"; Value = "  +0.37%  " }
    @{ Cell = "D45"; Value = "65.30" }
    @{ Cell = "E45"; Value = "  -2.79%  " }
    @{ Cell = "D46"; Value = "0.00000000119" }
    @{ Cell = "E46"; Value = "  +0.56%  " }
    @{ Cell = "D47"; Value = "7.069" }
    @{ Cell = "E47"; Value = "  -3.32%  " }
    @{ Cell = "D48"; Value = "0.3987" }
    @{ Cell = "E48"; Value = "  -0.62%  " }
    @{ Cell = "B49"; Value = "EnergySwap" }
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens" }
    @{ Cell = "D49"; Value = "8.898" }
    @{ Cell = "E49"; Value = "  -3.28%  " }
    @{ Cell = "B50"; Value = "RenderToken" }
    @{ Cell = "C50"; Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr" }
    @{ Cell = "D50"; Value = "1.666" }
    @{ Cell = "E50"; Value = "  -1.41%  " }
    @{ Cell = "E51"; Value = "  -0.26%  " }
)

foreach ($ch in $changes) {
    $rng = $ws.Range($ch.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $ch.Value
    $rng.Style = "Normal"
}
